$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Shows sheet: log a new show watched (La Palma)
# ---------------------------------------------------------------------
$wsShows = $wb.Worksheets.Item("Shows")
$wsShows.Activate()

$wsShows.Range("A19").Copy()
$wsShows.Range("A20").PasteSpecial(-4122) | Out-Null
$wsShows.Range("A20").Value = 45647
$wsShows.Range("B20").Value = "La Palma"
$wsShows.Range("C20").Value = 1
$wsShows.Range("D20").Value = 6.4
$wsShows.Range("E20").Value = "n"
$wsShows.Range("F20").Value = "Netflix"
$wsShows.Range("G20").Value = "Drama"
$wsShows.Range("H20").Value = "Meh. Kind of boring"

$wsShows.Range("A21").Select() | Out-Null

# ---------------------------------------------------------------------
# Movies sheet: log two new watched movies (Dec 21 + Dec 22)
# ---------------------------------------------------------------------
$wsMovies = $wb.Worksheets.Item("Movies")
$wsMovies.Activate()

$wsMovies.Range("A50").Copy()
$wsMovies.Range("A51").PasteSpecial(-4122) | Out-Null
$wsMovies.Range("A51").Value = 45647
$wsMovies.Range("B51").Value = "Ready or Not"
$wsMovies.Range("C51").Value = 7.8
$wsMovies.Range("D51").Value = "Horror/Thriller"
$wsMovies.Range("E51").Value = "y"
$wsMovies.Range("F51").Value = "netflix"
$wsMovies.Range("G51").Value = "good horror movie. And funny"

$wsMovies.Range("A50").Copy()
$wsMovies.Range("A52").PasteSpecial(-4122) | Out-Null
$wsMovies.Range("A52").Value = 45648
$wsMovies.Range("B52").Value = "I Spit on Your Grave"
$wsMovies.Range("C52").Value = 6.1
$wsMovies.Range("D52").Value = "Horror"
$wsMovies.Range("E52").Value = "y"
$wsMovies.Range("F52").Value = "netflix"

# ---------------------------------------------------------------------
# Poutine sheet: rating tweaks + a new vendor name
# ---------------------------------------------------------------------
$wsPoutine = $wb.Worksheets.Item("Poutine")
$wsPoutine.Activate()

$wsPoutine.Range("C4").Value = 8.7
$wsPoutine.Range("C6").Value = 8.4
$wsPoutine.Range("C10").Value = 8.9
$wsPoutine.Range("C12").Value = 9.2
$wsPoutine.Range("B14").Value = "Top Gun Burger"

$wsPoutine.Range("I22").Select() | Out-Null

# Back to Movies to fill in the last comment for the Dec 22 entry
$wsMovies.Activate()
$wsMovies.Range("G52").Value = "watched with denisse. Wow she is amazing at watching horrors this one was absolutely tough to watch. There was 2 rape scenes and at one of the major revenge scenes he cut off this guy's pen1s with garden scissors and shoved it in his mouth haha!"

$wsMovies.Range("G53").Select() | Out-Null

# ---------------------------------------------------------------------
# Books sheet: finish filling in "Play Nice" row + add "The Snakehead"
# ---------------------------------------------------------------------
$wsBooks = $wb.Worksheets.Item("Books")
$wsBooks.Activate()

$wsBooks.Range("C10").Copy()
$wsBooks.Range("C12").PasteSpecial(-4122) | Out-Null
$wsBooks.Range("C12").Value = 45587
$wsBooks.Range("E12").Value = "Jason Schreier"
$wsBooks.Range("F12").Value = "NF"
$wsBooks.Range("G12").Value = "y"
$wsBooks.Range("H12").Value = 8.5

$wsBooks.Range("D14").Value = "The Snakehead"
$wsBooks.Range("E14").Value = "Patrick Radden Keefe"
$wsBooks.Range("F14").Value = "NF"
$wsBooks.Range("G14").Value = "n"

$wsBooks.Range("I12").Value = "good book very interesting. Nice to hear the history of how everything rose and fall. Much lessons to be had"

$wsBooks.Range("D14").Select() | Out-Null

# ---------------------------------------------------------------------
# Golf sheet: just reposition the cursor, no data changes
# ---------------------------------------------------------------------
$wsGolf = $wb.Worksheets.Item("Golf")
$wsGolf.Activate()
$wsGolf.Range("H21").Select() | Out-Null

# ---------------------------------------------------------------------
# Finish back on Hockey (the tab that was active when the file was saved)
# ---------------------------------------------------------------------
$wsHockey = $wb.Worksheets.Item("Hockey")
$wsHockey.Activate()
$wsHockey.Range("H17").Select() | Out-Null
